# Update the 25 division problems in the practice-sheet table.
# Each data row of the table holds 5 problems; rows 1,5,9,13,17 (of 20)
# contain text, the others are blank spacer rows. One value ("878÷8=")
# occurs twice with two different replacements, so we address each
# cell individually by (row, column) rather than doing a single
# document-wide Find/Replace.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "827÷3="; New = "697÷5=" },
    @{ Row = 1;  Col = 2; Old = "402÷3="; New = "520÷2=" },
    @{ Row = 1;  Col = 3; Old = "152÷6="; New = "336÷8=" },
    @{ Row = 1;  Col = 4; Old = "376÷8="; New = "778÷7=" },
    @{ Row = 1;  Col = 5; Old = "974÷5="; New = "318÷6=" },

    @{ Row = 5;  Col = 1; Old = "515÷4="; New = "854÷9=" },
    @{ Row = 5;  Col = 2; Old = "262÷8="; New = "505÷6=" },
    @{ Row = 5;  Col = 3; Old = "878÷8="; New = "133÷7=" },
    @{ Row = 5;  Col = 4; Old = "891÷7="; New = "119÷5=" },
    @{ Row = 5;  Col = 5; Old = "790÷5="; New = "699÷6=" },

    @{ Row = 9;  Col = 1; Old = "132÷7="; New = "777÷7=" },
    @{ Row = 9;  Col = 2; Old = "243÷3="; New = "140÷3=" },
    @{ Row = 9;  Col = 3; Old = "433÷3="; New = "595÷7=" },
    @{ Row = 9;  Col = 4; Old = "198÷6="; New = "224÷4=" },
    @{ Row = 9;  Col = 5; Old = "221÷5="; New = "265÷5=" },

    @{ Row = 13; Col = 1; Old = "921÷8="; New = "623÷3=" },
    @{ Row = 13; Col = 2; Old = "844÷5="; New = "829÷4=" },
    @{ Row = 13; Col = 3; Old = "516÷4="; New = "145÷6=" },
    @{ Row = 13; Col = 4; Old = "184÷4="; New = "857÷9=" },
    @{ Row = 13; Col = 5; Old = "915÷6="; New = "831÷9=" },

    @{ Row = 17; Col = 1; Old = "991÷5="; New = "586÷8=" },
    @{ Row = 17; Col = 2; Old = "191÷4="; New = "772÷8=" },
    @{ Row = 17; Col = 3; Old = "316÷8="; New = "374÷2=" },
    @{ Row = 17; Col = 4; Old = "878÷8="; New = "790÷9=" },
    @{ Row = 17; Col = 5; Old = "305÷9="; New = "105÷8=" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cellRange = $cell.Range
    # Exclude the trailing end-of-cell marker so the new text takes the
    # place of exactly the old run's content, keeping formatting intact.
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    if ($textRange.Text -ne $item.Old) {
        throw "Unexpected cell text at row $($item.Row) col $($item.Col): [$($textRange.Text)] (expected [$($item.Old)])"
    }
    $textRange.Text = $item.New
}

Write-Output "Updated $($replacements.Count) problems"
